$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 4094803  # H33: 4679709.5 -> 4094803
$ws.Cells.Item(33, 9).Value = 4679685.5  # I33: 5459556.5 -> 4679685.5
$ws.Cells.Item(33, 11).Value = 4679685.5  # K33: 5459556.5 -> 4679685.5
$ws.Cells.Item(33, 13).Value = -4679456.5  # M33: -5459327.5 -> -4679456.5
$ws.Cells.Item(40, 8).Value = 7007.6  # H40: 7011 -> 7007.6
$ws.Cells.Item(40, 10).Value = 7370.75  # J40: 7496.3335 -> 7370.75
$ws.Cells.Item(40, 12).Value = 7370.75  # L40: 7496.3335 -> 7370.75
$ws.Cells.Item(40, 14).Value = -7720.75  # N40: -7846.3335 -> -7720.75
$ws.Cells.Item(43, 8).Value = 8738.333000000001  # H43: 8749.4 -> 8738.333000000001
$ws.Cells.Item(43, 10).Value = 8686.200000000001  # J43: 8687 -> 8686.200000000001
$ws.Cells.Item(43, 12).Value = 8686.200000000001  # L43: 8687 -> 8686.200000000001
$ws.Cells.Item(43, 14).Value = -8824.200000000001  # N43: -8825 -> -8824.200000000001
$ws.Cells.Item(62, 8).Value = 4962.125  # H62: 4974.875 -> 4962.125
$ws.Cells.Item(62, 9).Value = 4899  # I62: 4899.5 -> 4899
$ws.Cells.Item(62, 11).Value = 4899  # K62: 4899.5 -> 4899
$ws.Cells.Item(62, 13).Value = -4275  # M62: -4275.5 -> -4275
$ws.Cells.Item(65, 8).Value = 4962.125  # H65: 4974.875 -> 4962.125
$ws.Cells.Item(65, 9).Value = 4899  # I65: 4899.5 -> 4899
$ws.Cells.Item(65, 11).Value = 24495  # K65: 24497.5 -> 24495
$ws.Cells.Item(65, 13).Value = -21375  # M65: -21377.5 -> -21375
$ws.Cells.Item(98, 8).Value = 2405.4119  # H98: 2431 -> 2405.4119
$ws.Cells.Item(98, 9).Value = 2574.4167  # I98: 2627 -> 2574.4167
$ws.Cells.Item(98, 11).Value = 2574.4167  # K98: 2627 -> 2574.4167
$ws.Cells.Item(98, 13).Value = -1076.4167  # M98: -1129 -> -1076.4167
$ws.Cells.Item(122, 8).Value = 2405.4119  # H122: 2431 -> 2405.4119
$ws.Cells.Item(122, 9).Value = 2574.4167  # I122: 2627 -> 2574.4167
$ws.Cells.Item(122, 11).Value = 7723.250100000001  # K122: 7881 -> 7723.250100000001
$ws.Cells.Item(122, 13).Value = -5273.250100000001  # M122: -5431 -> -5273.250100000001
$ws.Cells.Item(132, 8).Value = 43300.2  # H132: 47057.305 -> 43300.2
$ws.Cells.Item(132, 9).Value = 50757.668  # I132: 55911.95 -> 50757.668
$ws.Cells.Item(132, 10).Value = 4148.5  # J132: 4997.75 -> 4148.5
$ws.Cells.Item(132, 11).Value = 152273.004  # K132: 167735.85 -> 152273.004
$ws.Cells.Item(132, 12).Value = 12445.5  # L132: 14993.25 -> 12445.5
$ws.Cells.Item(132, 13).Value = -149743.004  # M132: -165205.85 -> -149743.004
$ws.Cells.Item(132, 14).Value = -17505.5  # N132: -20053.25 -> -17505.5
$ws.Cells.Item(135, 8).Value = 761.2857  # H135: 761.3570999999999 -> 761.2857
$ws.Cells.Item(135, 9).Value = 918.2222  # I135: 918.3333 -> 918.2222
$ws.Cells.Item(135, 11).Value = 8263.9998  # K135: 8264.9997 -> 8263.9998
$ws.Cells.Item(135, 13).Value = -5728.9998  # M135: -5729.9997 -> -5728.9998

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2317  # H61: 2311.1 -> 2317
$ws.Cells.Item(61, 9).Value = 2317  # I61: 2314 -> 2317
$ws.Cells.Item(61, 10).Value = 0  # J61: 2285 -> 0
$ws.Cells.Item(61, 11).Value = 2317  # K61: 2314 -> 2317
$ws.Cells.Item(61, 12).Value = 0  # L61: 2285 -> 0
$ws.Cells.Item(61, 13).Value = -2105  # M61: -2102 -> -2105
$ws.Cells.Item(61, 14).ClearContents()  # N61: delete (was -2709)
$ws.Cells.Item(75, 8).Value = 0  # H75: 60173 -> 0
$ws.Cells.Item(75, 10).Value = 0  # J75: 60173 -> 0
$ws.Cells.Item(75, 12).Value = 0  # L75: 60173 -> 0
$ws.Cells.Item(75, 14).ClearContents()  # N75: delete (was -61921)
$ws.Cells.Item(76, 8).Value = 42632.668  # H76: 35579.8 -> 42632.668
$ws.Cells.Item(76, 10).Value = 42632.668  # J76: 35579.8 -> 42632.668
$ws.Cells.Item(76, 12).Value = 42632.668  # L76: 35579.8 -> 42632.668
$ws.Cells.Item(76, 14).Value = -43308.668  # N76: -36255.8 -> -43308.668
$ws.Cells.Item(78, 8).Value = 0  # H78: 60173 -> 0
$ws.Cells.Item(78, 10).Value = 0  # J78: 60173 -> 0
$ws.Cells.Item(78, 12).Value = 0  # L78: 180519 -> 0
$ws.Cells.Item(78, 14).ClearContents()  # N78: delete (was -189255)
$ws.Cells.Item(79, 8).Value = 42632.668  # H79: 35579.8 -> 42632.668
$ws.Cells.Item(79, 10).Value = 42632.668  # J79: 35579.8 -> 42632.668
$ws.Cells.Item(79, 12).Value = 42632.668  # L79: 35579.8 -> 42632.668
$ws.Cells.Item(79, 14).Value = -44972.668  # N79: -37919.8 -> -44972.668
$ws.Cells.Item(82, 8).Value = 60000  # H82: 0 -> 60000
$ws.Cells.Item(82, 10).Value = 60000  # J82: 0 -> 60000
$ws.Cells.Item(82, 12).Value = 60000  # L82: 0 -> 60000
$ws.Cells.Item(82, 14).Value = -60722  # N82: None -> -60722
$ws.Cells.Item(85, 8).Value = 60000  # H85: 0 -> 60000
$ws.Cells.Item(85, 10).Value = 60000  # J85: 0 -> 60000
$ws.Cells.Item(85, 12).Value = 60000  # L85: 0 -> 60000
$ws.Cells.Item(85, 14).Value = -62496  # N85: None -> -62496
$ws.Cells.Item(86, 8).Value = 30000  # H86: 0 -> 30000
$ws.Cells.Item(86, 10).Value = 30000  # J86: 0 -> 30000
$ws.Cells.Item(86, 12).Value = 30000  # L86: 0 -> 30000
$ws.Cells.Item(86, 14).Value = -32372  # N86: None -> -32372
$ws.Cells.Item(89, 8).Value = 30000  # H89: 0 -> 30000
$ws.Cells.Item(89, 10).Value = 30000  # J89: 0 -> 30000
$ws.Cells.Item(89, 12).Value = 90000  # L89: 0 -> 90000
$ws.Cells.Item(89, 14).Value = -101856  # N89: None -> -101856
$ws.Cells.Item(101, 8).Value = 0  # H101: 44900 -> 0
$ws.Cells.Item(101, 10).Value = 0  # J101: 44900 -> 0
$ws.Cells.Item(101, 12).Value = 0  # L101: 44900 -> 0
$ws.Cells.Item(101, 14).ClearContents()  # N101: delete (was -51390)
$ws.Cells.Item(102, 8).Value = 677.6  # H102: 622.7143 -> 677.6
$ws.Cells.Item(102, 9).Value = 677.6  # I102: 622.7143 -> 677.6
$ws.Cells.Item(102, 11).Value = 677.6  # K102: 622.7143 -> 677.6
$ws.Cells.Item(102, 13).Value = 944.4  # M102: 999.2857 -> 944.4
$ws.Cells.Item(132, 8).Value = 13517681  # H132: 13893160 -> 13517681
$ws.Cells.Item(132, 9).Value = 2647.926  # I132: 2732.8845 -> 2647.926
$ws.Cells.Item(132, 11).Value = 7943.778  # K132: 8198.6535 -> 7943.778
$ws.Cells.Item(132, 13).Value = -5413.778  # M132: -5668.6535 -> -5413.778
$ws.Cells.Item(136, 8).Value = 2317  # H136: 2311.1 -> 2317
$ws.Cells.Item(136, 9).Value = 2317  # I136: 2314 -> 2317
$ws.Cells.Item(136, 10).Value = 0  # J136: 2285 -> 0
$ws.Cells.Item(136, 11).Value = 6951  # K136: 6942 -> 6951
$ws.Cells.Item(136, 12).Value = 0  # L136: 6855 -> 0
$ws.Cells.Item(136, 13).Value = -4401  # M136: -4392 -> -4401
$ws.Cells.Item(136, 14).ClearContents()  # N136: delete (was -11955)

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 10002274  # H134: 10418974 -> 10002274
$ws.Cells.Item(134, 9).Value = 14707933  # I134: 15627086 -> 14707933
$ws.Cells.Item(134, 11).Value = 44123799  # K134: 46881258 -> 44123799
$ws.Cells.Item(134, 13).Value = -44121264  # M134: -46878723 -> -44121264

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3427.8823  # H31: 3343.5 -> 3427.8823
$ws.Cells.Item(34, 8).Value = 3427.8823  # H34: 3343.5 -> 3427.8823
$ws.Cells.Item(92, 8).Value = 96946.336  # H92: 98223.5 -> 96946.336
$ws.Cells.Item(92, 10).Value = 96946.336  # J92: 98223.5 -> 96946.336
$ws.Cells.Item(92, 12).Value = 96946.336  # L92: 98223.5 -> 96946.336
$ws.Cells.Item(92, 14).Value = -101938.336  # N92: -103215.5 -> -101938.336
$ws.Cells.Item(96, 8).Value = 47847.285  # H96: 49132.57 -> 47847.285
$ws.Cells.Item(96, 10).Value = 47847.285  # J96: 49132.57 -> 47847.285
$ws.Cells.Item(96, 12).Value = 47847.285  # L96: 49132.57 -> 47847.285
$ws.Cells.Item(96, 14).Value = -53339.285  # N96: -54624.57 -> -53339.285
$ws.Cells.Item(132, 8).Value = 3299.6  # H132: 3374.25 -> 3299.6
$ws.Cells.Item(132, 9).Value = 3000  # I132: 2999 -> 3000
$ws.Cells.Item(132, 11).Value = 9000  # K132: 8997 -> 9000
$ws.Cells.Item(132, 13).Value = -6470  # M132: -6467 -> -6470

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 108061.57  # H4: 100885.8 -> 108061.57
$ws.Cells.Item(4, 9).Value = 867.9048  # I4: 848.9091 -> 867.9048
$ws.Cells.Item(4, 10).Value = 429642.56  # J4: 375987.25 -> 429642.56
$ws.Cells.Item(4, 11).Value = 2603.7144  # K4: 2546.7273 -> 2603.7144
$ws.Cells.Item(4, 12).Value = 1288927.68  # L4: 1127961.75 -> 1288927.68
$ws.Cells.Item(4, 13).Value = -2491.7144  # M4: -2434.7273 -> -2491.7144
$ws.Cells.Item(4, 14).Value = -1289151.68  # N4: -1128185.75 -> -1289151.68
$ws.Cells.Item(26, 8).Value = 707.4286  # H26: 715.7143 -> 707.4286
$ws.Cells.Item(26, 9).Value = 594.3333  # I26: 400 -> 594.3333
$ws.Cells.Item(26, 10).Value = 761  # J26: 790 -> 761
$ws.Cells.Item(26, 11).Value = 1782.9999  # K26: 1200 -> 1782.9999
$ws.Cells.Item(26, 12).Value = 2283  # L26: 2370 -> 2283
$ws.Cells.Item(26, 13).Value = -1494.9999  # M26: -912 -> -1494.9999
$ws.Cells.Item(26, 14).Value = -2859  # N26: -2946 -> -2859
$ws.Cells.Item(70, 8).Value = 1673.125  # H70: 1899.8334 -> 1673.125
$ws.Cells.Item(70, 9).Value = 1683.5714  # I70: 1899.8334 -> 1683.5714
$ws.Cells.Item(70, 10).Value = 1600  # J70: 0 -> 1600
$ws.Cells.Item(70, 11).Value = 5050.7142  # K70: 5699.5002 -> 5050.7142
$ws.Cells.Item(70, 12).Value = 4800  # L70: 0 -> 4800
$ws.Cells.Item(70, 13).Value = -4735.7142  # M70: -5384.5002 -> -4735.7142
$ws.Cells.Item(70, 14).Value = -5430  # N70: None -> -5430
$ws.Cells.Item(73, 8).Value = 1673.125  # H73: 1899.8334 -> 1673.125
$ws.Cells.Item(73, 9).Value = 1683.5714  # I73: 1899.8334 -> 1683.5714
$ws.Cells.Item(73, 10).Value = 1600  # J73: 0 -> 1600
$ws.Cells.Item(73, 11).Value = 5050.7142  # K73: 5699.5002 -> 5050.7142
$ws.Cells.Item(73, 12).Value = 4800  # L73: 0 -> 4800
$ws.Cells.Item(73, 13).Value = -3958.7142  # M73: -4607.5002 -> -3958.7142
$ws.Cells.Item(73, 14).Value = -6984  # N73: None -> -6984

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(39, 8).Value = 4999.5  # H39: 5000 -> 4999.5
$ws.Cells.Item(39, 10).Value = 4999.5  # J39: 5000 -> 4999.5
$ws.Cells.Item(39, 12).Value = 4999.5  # L39: 5000 -> 4999.5
$ws.Cells.Item(39, 14).Value = -6063.5  # N39: -6064 -> -6063.5
$ws.Cells.Item(132, 8).Value = 2018.6316  # H132: 2214.0667 -> 2018.6316
$ws.Cells.Item(132, 9).Value = 1829.9166  # I132: 1962.2222 -> 1829.9166
$ws.Cells.Item(132, 10).Value = 2342.1428  # J132: 2591.8333 -> 2342.1428
$ws.Cells.Item(132, 11).Value = 5489.7498  # K132: 5886.6666 -> 5489.7498
$ws.Cells.Item(132, 12).Value = 7026.428400000001  # L132: 7775.499899999999 -> 7026.428400000001
$ws.Cells.Item(132, 13).Value = -2959.7498  # M132: -3356.6666 -> -2959.7498
$ws.Cells.Item(132, 14).Value = -12086.4284  # N132: -12835.4999 -> -12086.4284

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(32, 8).Value = 11579.429  # H32: 11342.833 -> 11579.429
$ws.Cells.Item(32, 9).Value = 11666.333  # I32: 11000 -> 11666.333
$ws.Cells.Item(32, 11).Value = 11666.333  # K32: 11000 -> 11666.333
$ws.Cells.Item(32, 13).Value = -11349.333  # M32: -10683 -> -11349.333
$ws.Cells.Item(34, 8).Value = 22982.334  # H34: 27299.666 -> 22982.334
$ws.Cells.Item(34, 9).Value = 25474.75  # I34: 27299.666 -> 25474.75
$ws.Cells.Item(34, 10).Value = 17997.5  # J34: 0 -> 17997.5
$ws.Cells.Item(34, 11).Value = 25474.75  # K34: 27299.666 -> 25474.75
$ws.Cells.Item(34, 12).Value = 17997.5  # L34: 0 -> 17997.5
$ws.Cells.Item(34, 13).Value = -25271.75  # M34: -27096.666 -> -25271.75
$ws.Cells.Item(34, 14).Value = -18403.5  # N34: None -> -18403.5
$ws.Cells.Item(96, 8).Value = 1921.75  # H96: 1933 -> 1921.75
$ws.Cells.Item(96, 10).Value = 1896  # J96: 1900 -> 1896
$ws.Cells.Item(96, 12).Value = 1896  # L96: 1900 -> 1896
$ws.Cells.Item(96, 14).Value = -4642  # N96: -4646 -> -4642
$ws.Cells.Item(103, 8).Value = 25000  # H103: 19999.5 -> 25000
$ws.Cells.Item(103, 10).Value = 25000  # J103: 19999.5 -> 25000
$ws.Cells.Item(103, 12).Value = 25000  # L103: 19999.5 -> 25000
$ws.Cells.Item(103, 14).Value = -27344  # N103: -22343.5 -> -27344
$ws.Cells.Item(122, 8).Value = 1880.7826  # H122: 1966.5 -> 1880.7826
$ws.Cells.Item(122, 9).Value = 1986.6666  # I122: 2103.8235 -> 1986.6666
$ws.Cells.Item(122, 11).Value = 5959.9998  # K122: 6311.470499999999 -> 5959.9998
$ws.Cells.Item(122, 13).Value = -3509.9998  # M122: -3861.470499999999 -> -3509.9998
